$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "sep18" sheet by duplicating "aug18" (so it inherits
#    the exact same formatting: column width/bestFit, date style, cell
#    styles, formulas, page margins, etc.), placed right after "aug18".
# ---------------------------------------------------------------------
$aug = $wb.Worksheets.Item("aug18")
$aug.Copy([System.Reflection.Missing]::Value, $aug)
$sep = $wb.Worksheets.Item($aug.Index + 1)
$sep.Name = "sep18"

# ---------------------------------------------------------------------
# 2. Update week 1 (row 1 header date + rows 2-19 scores) for sep18.
# ---------------------------------------------------------------------
$sep.Range("A1").Value = 43352

$week1 = @(
    @(5, "S", 3),
    @(3, "", 2),
    @(4, "R", 1),
    @(5, "R", 1),
    @(4, "", 2),
    @(5, "S", 2),
    @(4, "S", 2),
    @(3, "", 1),
    @(5, "S", 2),
    @(5, "S", 2),
    @(5, "", 1),
    @(5, "S", 2),
    @(5, "L", 1),
    @(7, "R", 2),
    @(5, "", 1),
    @(4, "R", 2),
    @(4, "R", 2),
    @(5, "S", 2)
)

$r = 2
foreach ($row in $week1) {
    $sep.Cells.Item($r, 2).Value = $row[0]
    if ($row[1] -ne "") {
        $sep.Cells.Item($r, 3).Value = $row[1]
    } else {
        $sep.Cells.Item($r, 3).ClearContents()
    }
    $sep.Cells.Item($r, 5).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------
# 3. Update week 2 (row 22 header date + rows 23-40 scores) for sep18.
# ---------------------------------------------------------------------
$sep.Range("A22").Value = 43359

$week2 = @(
    @(5, "S", 3),
    @(3, "", 1),
    @(5, "S", 3),
    @(5, "R", 1),
    @(3, "", 1),
    @(7, "L", 3),
    @(5, "R", 2),
    @(7, "", 2),
    @(5, "S", 2),
    @(5, "L", 3),
    @(3, "", 1),
    @(5, "S", 2),
    @(4, "S", 2),
    @(7, "S", 2),
    @(3, "", 1),
    @(3, "R", 1),
    @(4, "S", 2),
    @(4, "S", 1)
)

$r = 23
foreach ($row in $week2) {
    $sep.Cells.Item($r, 2).Value = $row[0]
    if ($row[1] -ne "") {
        $sep.Cells.Item($r, 3).Value = $row[1]
    } else {
        $sep.Cells.Item($r, 3).ClearContents()
    }
    $sep.Cells.Item($r, 5).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------
# 4. Week 3 (rows 43-61) hasn't happened yet - clear the date and all
#    scores, leaving just the header labels and the hole names.
# ---------------------------------------------------------------------
$sep.Range("A43").ClearContents()
$sep.Range("B44:B61").ClearContents()
$sep.Range("C44:C61").ClearContents()
$sep.Range("E44:E61").ClearContents()

# ---------------------------------------------------------------------
# 5. Selection/view state on the new sheet.
# ---------------------------------------------------------------------
$sep.Range("G38").Select()

# ---------------------------------------------------------------------
# 6. Update the "aug18" sheet's view (no longer the active tab).
# ---------------------------------------------------------------------
$aug.Activate()
$excel.ActiveWindow.ScrollRow = 13
$aug.Range("C67").Select()

# ---------------------------------------------------------------------
# 7. Update the "dec 17" sheet's view (no longer tabSelected either).
# ---------------------------------------------------------------------
$dec17 = $wb.Worksheets.Item("dec 17")
$dec17.Activate()
$dec17.Range("G28").Select()

# ---------------------------------------------------------------------
# 8. Leave "sep18" as the final active sheet/tab.
# ---------------------------------------------------------------------
$sep.Activate()
$sep.Range("G38").Select()
